$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two rows: A1 holds a numeric placeholder styled with
# a bold font + thin border (cellXf index 1), and A2 holds the shared string
# with the question payload (plain default style). The edit removes the old
# A1 row completely - Range.Delete() shifts A2 up into A1, which drops the
# placeholder value *and* its special font/border style in one step, leaving
# A1 with the default (unstyled) format - then the payload text is rewritten
# as pretty-printed JSON (still single-quoted Python dict syntax replaced by
# double-quoted, indented JSON; same question data).
$ws.Range("A1").Delete()

$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are a product marketing manager for an industrial software company. You and your team have recently created well-rounded customer personas, which should be used to improve the companywide understanding of your ideal customer base. However, you have not yet shared them with the rest of the company.Which action should you take when doing this?",
        "ques_type": 2,
        "options": [
            "Add a document detailing the personas to your company\u2019s internal database. ",
            "Host tailored training sessions on the personas for each department. ",
            "Email every stakeholder to announce the creation and details of the personas.",
            "Distribute printed handouts that detail the personas companywide."
        ],
        "score": "Host tailored training sessions on the personas for each department."
    },
    {
        "title": "You work for a B2B software-as-a-service (SaaS) company that provides financial management software. You need to evaluate how effectively your users employ the software and identify areas where additional education is required.Which actions should you take?",
        "ques_type": 15,
        "options": [
            "Identify usage patterns using web analytics.",
            "Gather input from social media.",
            "Monitor user activity on your analytics dashboard.",
            "Enhance the software\u2019s capabilities based on well-accepted user experience frameworks. ",
            "Conduct user surveys."
        ],
        "score": [
            "Identify usage patterns using web analytics.",
            "Conduct user surveys."
        ]
    },
    {
        "title": "You are a product marketing manager for a B2B employee engagement software. Your sales team is struggling to convert leads to customers, so you want to create case studies to emphasize your product''s value. For this, you need to work with the customer success team. Which essential step should you take?",
        "ques_type": 2,
        "options": [
            "Identify successful customers and their success metrics.",
            "Assign each case study to a customer success manager. ",
            "Create each case study based on your knowledge. ",
            "Update the terms and conditions of the product. "
        ],
        "score": "Identify successful customers and their success metrics."
    },
    {
        "title": "You are a product marketing manager for a B2B software-as-a-service (SaaS) tool with a big sales team. You have been running weekly sales training sessions, but attendance has been dropping even though you still have content to deliver. Last week, you sent out an anonymous survey, and the results suggested that the sessions had become \u201cboring,\u201d \u201cof inconsistent relevance,\u201d and \u201crepetitive.\u201dWhich action should you take?",
        "ques_type": 2,
        "options": [
            "Switch to a monthly cadence. ",
            "Postpone the sales training until interest increases again.",
            "Recruit a sales colleague to run the next session with you. ",
            "Make attendance for the meeting mandatory."
        ],
        "score": "Recruit a sales colleague to run the next session with you."
    }
]'

# Re-run autofit so the row height stays at the sheet default instead of
# growing to fit the many embedded newlines in the new multi-line text.
$ws.Rows.Item(1).AutoFit()
